# Append: 2025-10-25 01:42 JST
# - Refresh the "取得日時" (fetched-at) timestamp on every existing data row.
# - Insert a brand-new scraped listing as the new row 16 ("運用中HPのドメイン分け"),
#   which pushes the previous rows 16-18 down to 17-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newStamp = "2025-10-25 01:42:45"

# 1) Refresh timestamps for every existing data row (rows 2-18) before the
#    insert shifts anything, so every row ends up carrying the new stamp.
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = $newStamp
}

# 2) Insert a new row above row 16, shifting the old rows 16-18 down to 17-19.
$ws.Rows.Item(16).Insert()

# 3) Populate the newly inserted row 16 with the new listing. (Only the URL
#    *text* lands on row 16 here - the matching Hyperlink relationship is
#    registered down on row 19, mirroring the source data exactly.)
$ws.Cells.Item(16, 1).Value = $newStamp
$ws.Cells.Item(16, 2).Value = "運用中HPのドメイン分け"
$ws.Cells.Item(16, 3).Value = "システム開発"
$ws.Cells.Item(16, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(16, 5).Value = "期限情報なし"
$newUrl = "https://www.lancers.jp/work/detail/5420277"
$ws.Cells.Item(16, 6).Value = $newUrl
$ws.Cells.Item(16, 7).Value = 13

# 4) Register the hyperlink relationship for the new listing.
$ws.Hyperlinks.Add($ws.Cells.Item(19, 6), $newUrl)
$ws.Cells.Item(19, 6).Style = "Hyperlink"

$wb.Save()
